$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of user data (row 8)
$ws.Range("A8").Value = "316904978"
$ws.Range("A8").NumberFormat = "@"
$ws.Range("B8").Value = "Michael"
$ws.Range("C8").Value = "Elisha"
$ws.Range("D8").Value = "123abc1a"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").Value = $true

# Update selection to reflect the active cell used when saving
[void]$ws.Range("D5").Select()
